$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column cells hold numeric-looking text (e.g. "28.428.44",
# "0.0583") that must stay plain text, matching the original inline-string
# cell type, instead of being auto-coerced to a number by Excel. Forcing
# the NumberFormat to Text ("@") before the assignment prevents that
# coercion; ClearFormats() afterwards drops the now-unneeded format so the
# cell keeps its original (unstyled) appearance.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.428.44'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.549.52'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.15%  '

$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.40'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.481'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.00%  '

$ws.Range("E7").Value = '  -0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.95'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.38%  '

$ws.Range("E9").Value = '  -1.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0583'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0888'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.770.33'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.554.62'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.94%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.411.13'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.18%  '

$ws.Range("E15").Value = '  -2.14%  '

$ws.Range("E16").Value = '  -2.37%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.94'
$ws.Range("D17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.23'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.13%  '

$ws.Range("E19").Value = '  -1.79%  '

$ws.Range("E20").Value = '  -2.74%  '

$ws.Range("E21").Value = '  -0.27%  '

$ws.Range("E22").Value = '  -0.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.90'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.71%  '

$ws.Range("E24").Value = '  -2.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.47'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.72'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.103'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.51%  '

$ws.Range("E28").Value = '  -0.24%  '

$ws.Range("E29").Value = '  -3.78%  '

$ws.Range("E31").Value = '  -4.41%  '

$ws.Range("E32").Value = '  -1.96%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.382.98'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.17%  '

$ws.Range("E34").Value = '  -3.68%  '

$ws.Range("E35").Value = '  -1.64%  '

$ws.Range("E36").Value = '  -3.34%  '

$ws.Range("E37").Value = '  -2.77%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.57'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.40%  '

$ws.Range("E39").Value = '  -2.57%  '

$ws.Range("E40").Value = '  +1.69%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.769'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.69%  '

$ws.Range("E44").Value = '  -1.56%  '

$ws.Range("E45").Value = '  -2.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.71'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.683.41'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.872'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -9.16%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.18'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '42.88'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +7.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0102'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.01%  '
